# Jogos_do_Dia_Betfair_Back_Lay_2025-12-24.xlsx update
# - refresh odds for the three existing Thai League 1 matches (rows 2-4)
# - append two new Saudi 1st Division matches (rows 5-6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Ratchaburi x Chiangrai Utd) odds refresh ---
$ws.Range("F2").Value  = 1.56
$ws.Range("G2").Value  = 1.66
$ws.Range("H2").Value  = 5.6
$ws.Range("I2").Value  = 6.8
$ws.Range("K2").Value  = 5
$ws.Range("N2").Value  = 4
$ws.Range("O2").Value  = 1.26
$ws.Range("P2").Value  = 2.06
$ws.Range("Q2").Value  = 1.77
$ws.Range("R2").Value  = 1.41
$ws.Range("S2").Value  = 3
$ws.Range("U2").Value  = 1.92
$ws.Range("V2").Value  = 1.17
$ws.Range("W2").Value  = 2.5
$ws.Range("AN2").Value = 10

# --- Row 3 (Sukhothai x Buriram Utd) odds refresh ---
$ws.Range("G3").Value  = 10
$ws.Range("I3").Value  = 1.45
$ws.Range("J3").Value  = 5.1
$ws.Range("K3").Value  = 6
$ws.Range("L3").Value  = 1.24
$ws.Range("N3").Value  = 5.1
$ws.Range("P3").Value  = 2.42
$ws.Range("Q3").Value  = 1.6
$ws.Range("R3").Value  = 1.56
$ws.Range("T3").Value  = 1.87
$ws.Range("V3").Value  = 3.2
$ws.Range("AB3").Value = 40
$ws.Range("AD3").Value = 10.5
$ws.Range("AG3").Value = 40

# --- Row 4 (BG Pathumthani United x Dragon Pathumwan Kanchana) odds refresh ---
$ws.Range("G4").Value = 44
$ws.Range("H4").Value = 1.1
$ws.Range("S4").Value = 1.66
$ws.Range("T4").Value = 1.83
$ws.Range("V4").Value = 1.02
$ws.Range("W4").Value = 1.02

# --- Row 5 (new): Saudi 1st Division - Al-Anwar Club x Al Batin ---
$ws.Range("A5").Value = "Saudi 1st Division"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2025-12-24"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "09:10:00"
$ws.Range("D5").Value = "Al-Anwar Club"
$ws.Range("E5").Value = "Al Batin"
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 1.04
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 1.04
$ws.Range("K5").Value = 950
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 1.16
$ws.Range("O5").Value = 1.01
$ws.Range("P5").Value = 1.15
$ws.Range("Q5").Value = 1.01
$ws.Range("R5").Value = 1.08
$ws.Range("S5").Value = 1.02
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 1.01
$ws.Range("V5").Value = 1.02
$ws.Range("W5").Value = 1.02
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# --- Row 6 (new): Saudi 1st Division - Al Bukayriyah x Al-Arabi Al-Saudi ---
$ws.Range("A6").Value = "Saudi 1st Division"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2025-12-24"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "09:20:00"
$ws.Range("D6").Value = "Al Bukayriyah"
$ws.Range("E6").Value = "Al-Arabi Al-Saudi"
$ws.Range("F6").Value = 1.74
$ws.Range("G6").Value = 2.18
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 6.4
$ws.Range("J6").Value = 3.25
$ws.Range("K6").Value = 3.8
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 1.54
$ws.Range("O6").Value = 1.01
$ws.Range("P6").Value = 1.54
$ws.Range("Q6").Value = 2.04
$ws.Range("R6").Value = 1.18
$ws.Range("S6").Value = 2.04
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.18
$ws.Range("W6").Value = 1.85
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000
